# Edit script: apply the changes described by the diff to the active document.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Update the date mentioned in the introductory paragraph (19 -> 24 de mayo de 2022)
Replace-Text "suscrito con fecha  19 de mayo de 2022" "suscrito con fecha  24 de mayo de 2022"

# 2. Update the second data row first (to avoid collisions with row 1's new values)
Replace-Text "Moreno Ramos Laura" "Díez Viñas Malena"
Replace-Text "13c" "14d"
# "fd" must be replaced before "segdzg" because the replacement text "sdfdsf"
# itself contains the substring "fd", which would otherwise get re-matched.
Replace-Text "fd" "dsfds"
Replace-Text "segdzg" "sdfdsf"
Replace-Text "2022-03-04" "2022-05-24"
Replace-Text "2022-03-05" "2022-05-26"

# 3. Update the first data row
Replace-Text "Jiménez Coello Daniel" "Moreno Ramos Laura"
Replace-Text "11a" "13c"
Replace-Text "2022-05-03" "2022-05-24"
Replace-Text "2022-05-11" "2022-05-26"

# 4. Update the closing date line ("En Puertollano a  19  de mayo  2022")
Replace-Text "En Puertollano a  19  de mayo  2022" "En Puertollano a  24  de mayo  2022"
